$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.575.19"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "2.139.02"
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "'352.27"
$ws.Range("E5").Value = "  +5.42%  "
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").Value = "'0.5262"
$ws.Range("E7").Value = "  +0.95%  "
$ws.Range("D8").Value = "'0.4572"
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").Value = "'53.49"
$ws.Range("E9").Value = "  -1.95%  "
$ws.Range("D10").Value = "'0.09159"
$ws.Range("E10").Value = "  +3.02%  "
$ws.Range("D11").Value = "'1.189"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "'25.51"
$ws.Range("E12").Value = "  +5.96%  "
$ws.Range("D13").Value = "2.135.63"
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("D14").Value = "'6.903"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").Value = "'8.192"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").Value = "'102.30"
$ws.Range("E16").Value = "  +5.40%  "
$ws.Range("E17").Value = "  +2.61%  "
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "'0.06729"
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").Value = "'20.55"
$ws.Range("E20").Value = "  +7.18%  "
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("D23").Value = "30.651.10"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").Value = "'12.93"
$ws.Range("E24").Value = "  +4.90%  "
$ws.Range("E25").Value = "  +2.30%  "
$ws.Range("D26").Value = "2.390.78"
$ws.Range("E26").Value = "  +2.12%  "
$ws.Range("D27").Value = "'22.59"
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("D28").Value = "'2.648"
$ws.Range("E28").Value = "  +5.46%  "
$ws.Range("D29").Value = "'164.96"
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("D30").Value = "'136.09"
$ws.Range("E30").Value = "  +2.29%  "
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("D32").Value = "'0.1084"
$ws.Range("E32").Value = "  +1.77%  "
$ws.Range("D33").Value = "'1.691"
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("D34").Value = "'6.429"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").Value = "'4.043"
$ws.Range("E35").Value = "  +2.51%  "
$ws.Range("D36").Value = "'6.125"
$ws.Range("E36").Value = "  +5.80%  "
$ws.Range("D37").Value = "'10.53"
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("D38").Value = "'0.02651"
$ws.Range("E38").Value = "  +3.05%  "
$ws.Range("D39").Value = "'0.06985"
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("D40").Value = "'0.2341"
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("D41").Value = "'12.76"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").Value = "'0.7015"
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("D43").Value = "'1.281"
$ws.Range("E43").Value = "  +2.92%  "
$ws.Range("D44").Value = "'14.82"
$ws.Range("E44").Value = "  +6.20%  "
$ws.Range("D45").Value = "'2.368"
$ws.Range("E45").Value = "  +2.23%  "
$ws.Range("D46").Value = "'0.6528"
$ws.Range("E46").Value = "  +2.89%  "
$ws.Range("D47").Value = "'0.00000000371"
$ws.Range("E47").Value = "  +6.92%  "
$ws.Range("D48").Value = "'3.757"
$ws.Range("E48").Value = "  +2.88%  "
$ws.Range("D49").Value = "'1.252"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("D50").Value = "'84.11"
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("D51").Value = "'0.07301"
$ws.Range("E51").Value = "  +2.58%  "
